# Update isSourceOf and isDerivedFrom headings to correct (capitalized) form.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: C1 was "isDerivedFrom" -> "IsDerivedFrom"
#             D1 was "isSourceOf"   -> "IsSourceOf"
$ws.Range("C1").Value = "IsDerivedFrom"
$ws.Range("D1").Value = "IsSourceOf"

# Update the active selection to D2 (matches the saved selection in the file)
$ws.Range("D2").Select()
